$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1.xml)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2900.628494008129
$ws.Range("E2").Value = 290490.7128553893
$ws.Range("G2").Value = 80959.25712661196
$ws.Range("I2").Value = 149420.986357725
$ws.Range("L2").Value = 509988.6069102
$ws.Range("M2").Value = 112287.0813999
$ws.Range("N2").Value = 71616.34392528806
$ws.Range("O2").Value = 66869.92067293868

# Sheet "2030" (sheet2.xml)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 2297.730639432591
$ws.Range("B2").Value = 35125.66274683856
$ws.Range("E2").Value = 164480.3067901611
$ws.Range("I2").Value = 162514.5034018797
$ws.Range("L2").Value = 92628.68888285091
$ws.Range("M2").Value = 61433.01601085002
$ws.Range("N2").Value = 19245.54644840866
$ws.Range("O2").Value = 11607.72427391814

# Sheet "2035" (sheet3.xml)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 9004.007062543757
$ws.Range("B2").Value = 26188.96215853551
$ws.Range("E2").Value = 139750.6423037671
$ws.Range("I2").Value = 167681.5062414853
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 66005.4619483283
$ws.Range("N2").Value = 43781.27610165381
$ws.Range("O2").Value = 52207.37179957087

# Sheet "2040" (sheet4.xml)
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 243.8709858932979

# Sheet "2045" (sheet5.xml)
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 35564.54260491626
$ws.Range("N2").Value = 0
